# "Generate Report for Handoff"
#
# The ffb65260-... file has finished translation and is now ready for
# handoff. The status report rows for the three files that share this
# handoff batch (8db96cd7-..., c93166b9-..., ffb65260-...) are refreshed:
# all three now show "Ready for handoff" and the row order / handoff
# timestamps are updated to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A7").Value = "8db96cd7-a28c-45e8-9a11-a35e610ed50b.md"
$overview.Range("B7").Value = "Ready for handoff"
$overview.Range("C7").Value = "Ready for handoff"
$overview.Range("D7").Value = "2016-32-12 04:32:59"

$overview.Range("A8").Value = "c93166b9-b152-4ed2-9188-d7e4c736fc46.md"
$overview.Range("B8").Value = "Ready for handoff"
$overview.Range("C8").Value = "Ready for handoff"
$overview.Range("D8").Value = "2016-31-12 04:31:15"

$overview.Range("A9").Value = "ffb65260-b2ed-415b-9a80-58322094462c.md"
$overview.Range("B9").Value = "Ready for handoff"
$overview.Range("C9").Value = "Ready for handoff"
$overview.Range("D9").Value = "2016-35-12 04:35:28"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A7").Value = "8db96cd7-a28c-45e8-9a11-a35e610ed50b.md"
$zhcn.Range("C7").Value = "Ready for handoff"
$zhcn.Range("D7").Value = "8db96cd7-a28c-45e8-9a11-a35e610ed50b.fcc14e03234585534efe5559461e58200afcf7d8.zh-cn.xlf"
$zhcn.Range("E7").Value = "2016-03-12 04:32:56"

$zhcn.Range("A8").Value = "c93166b9-b152-4ed2-9188-d7e4c736fc46.md"
$zhcn.Range("C8").Value = "Ready for handoff"
$zhcn.Range("D8").Value = "c93166b9-b152-4ed2-9188-d7e4c736fc46.608a58dbbce996c93cec27acc58bd782e7ac473d.zh-cn.xlf"
$zhcn.Range("E8").Value = "2016-03-12 04:31:12"

$zhcn.Range("A9").Value = "ffb65260-b2ed-415b-9a80-58322094462c.md"
$zhcn.Range("C9").Value = "Ready for handoff"
$zhcn.Range("D9").Value = "ffb65260-b2ed-415b-9a80-58322094462c.4e9b268e2c3f9cbc74702e04cdc12db8f3b78063.zh-cn.xlf"
$zhcn.Range("E9").Value = "2016-03-12 04:35:25"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A7").Value = "8db96cd7-a28c-45e8-9a11-a35e610ed50b.md"
$dede.Range("C7").Value = "Ready for handoff"
$dede.Range("D7").Value = "8db96cd7-a28c-45e8-9a11-a35e610ed50b.fcc14e03234585534efe5559461e58200afcf7d8.de-de.xlf"
$dede.Range("E7").Value = "2016-03-12 04:32:59"

$dede.Range("A8").Value = "c93166b9-b152-4ed2-9188-d7e4c736fc46.md"
$dede.Range("C8").Value = "Ready for handoff"
$dede.Range("D8").Value = "c93166b9-b152-4ed2-9188-d7e4c736fc46.608a58dbbce996c93cec27acc58bd782e7ac473d.de-de.xlf"
$dede.Range("E8").Value = "2016-03-12 04:31:15"

$dede.Range("A9").Value = "ffb65260-b2ed-415b-9a80-58322094462c.md"
$dede.Range("C9").Value = "Ready for handoff"
$dede.Range("D9").Value = "ffb65260-b2ed-415b-9a80-58322094462c.4e9b268e2c3f9cbc74702e04cdc12db8f3b78063.de-de.xlf"
$dede.Range("E9").Value = "2016-03-12 04:35:28"
